# Fruta / hortaliza, semanal
# Insert a new data row (new row 4) shifting the existing rows 4-10 down to 5-11,
# and populate the new row with the latest weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4; this pushes current rows 4-10 down to 5-11
# and Excel carries down the row-above formatting (e.g. the date style on column D).
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new weekly record.
$ws.Range("A4").Value = 10
$ws.Range("B4").Value = "Vega Modelo de Temuco"
$ws.Range("C4").Value = "La Araucanía"
$ws.Range("D4").Value = 44616
$ws.Range("E4").Value = 9
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100101
$ws.Range("H4").Value = "Berries"
$ws.Range("I4").Value = 100101004
$ws.Range("J4").Value = "Frambuesa"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 3200
$ws.Range("O4").Value = 3200
$ws.Range("P4").Value = 3200
$ws.Range("Q4").Value = "$/envase 1 kilo"
$ws.Range("R4").Value = "Región de La Araucanía"
$ws.Range("S4").Value = 3200
$ws.Range("T4").Value = 1
